$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2024-09-22 23:48:05", "check_availability", "https://www.opentable.com/r/the-rux-nashville", "Checked availability: No availability for the selected date.", "2024-09-22", "23:48:05"),
    @("2024-09-22 23:48:58", "check_availability", "https://www.opentable.com/r/the-rux-nashville", "Checked availability: No availability for the selected date.", "2024-09-22", "23:48:58"),
    @("2024-09-22 23:49:21", "check_availability", "https://www.opentable.com/r/the-rux-nashville", "Checked availability: No availability for the selected date.", "2024-09-22", "23:49:21")
)

$startRow = 6
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Column E holds plain "yyyy-mm-dd" text that Excel would otherwise
        # auto-convert to a date serial number; force it to stay text, then
        # reset the cell style back to Normal so no extra formatting sticks.
        if ($c -eq 5) {
            $cell.NumberFormat = "@"
            $cell.Value = $rowData[$c - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $rowData[$c - 1]
        }
    }
}
